# Weekly update: a new price record was added for "Perejil" (Terminal La
# Palmera de La Serena) and all the subsequent weekly records shift down
# by one row (row 28 becomes the new entry; former rows 28..121 become
# 29..122; dimension grows from A1:R121 to A1:R122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 28 - this shifts rows 28:121
# down to 29:122 and (as in Excel) the new row inherits the formatting
# (date style) of the row that used to be above it.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with this week's record.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44607
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112044
$ws.Range("G28").Value = "Perejil"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 2400
$ws.Range("K28").Value = 2300
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2400
$ws.Range("N28").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O28").Value = "Provincia del Elquí"
$ws.Range("P28").Value = 1600
$ws.Range("Q28").Value = 1.5
$ws.Range("R28").Value = "Hortaliza"
